$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "BIRTHPLACE, CITY"
$ws.Range("D3").Value = "this table consists of patient data admittied in clinics, using cte, can you check which BIRTHPLACE are matching with CITY."

$ws.Range("C3").Select()
